# Auto-generated edit script for contratos-7-2019.xlsx
# Fixes:
#  1) Replace commas with periods in 4 'Razon social' entries (column E)
#  2) Reformat 263 Argentine-formatted amount strings in column H
#     ('1.234,56' -> '1234.56') while keeping them as TEXT shared-strings
#     (no number-format/style change), matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: punctuation fixes in 'Razon social' (column E) ---
$ws.Range("E53").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E112").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E85").Value = "URUMAT SOCIEDAD SIMPLE DE BONASEGLA CATALINA. BONASEGLA LUCIANA Y BONASEGLA SILVIO"
$ws.Range("E217").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E236").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E219").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

# --- Part 2: numeric-string reformatting in 'Importe' (column H) ---
# Values are entered with a leading apostrophe so Excel keeps them as text
# (otherwise '122400.00' would be auto-converted to the number 122400 and
# lose its trailing '.00'); ClearFormats() then strips the resulting
# quote-prefix text styling so the cell keeps the sheet's default style,
# exactly like the original cells (which also carry no explicit style).
$ws.Range("H2").Value = "'122400.00"
$ws.Range("H2").ClearFormats()
$ws.Range("H4").Value = "'5100.00"
$ws.Range("H4").ClearFormats()
$ws.Range("H5").Value = "'7200.00"
$ws.Range("H5").ClearFormats()
$ws.Range("H6").Value = "'46835.00"
$ws.Range("H6").ClearFormats()
$ws.Range("H7").Value = "'47949.92"
$ws.Range("H7").ClearFormats()
$ws.Range("H8").Value = "'316000.00"
$ws.Range("H8").ClearFormats()
$ws.Range("H9").Value = "'772000.00"
$ws.Range("H9").ClearFormats()
$ws.Range("H10").Value = "'3999.00"
$ws.Range("H10").ClearFormats()
$ws.Range("H11").Value = "'220.00"
$ws.Range("H11").ClearFormats()
$ws.Range("H12").Value = "'3560.00"
$ws.Range("H12").ClearFormats()
$ws.Range("H13").Value = "'1759.00"
$ws.Range("H13").ClearFormats()
$ws.Range("H14").Value = "'10350.00"
$ws.Range("H14").ClearFormats()
$ws.Range("H15").Value = "'709500.00"
$ws.Range("H15").ClearFormats()
$ws.Range("H16").Value = "'1116.41"
$ws.Range("H16").ClearFormats()
$ws.Range("H17").Value = "'796250.00"
$ws.Range("H17").ClearFormats()
$ws.Range("H18").Value = "'654561.60"
$ws.Range("H18").ClearFormats()
$ws.Range("H21").Value = "'605000.00"
$ws.Range("H21").ClearFormats()
$ws.Range("H22").Value = "'3151.00"
$ws.Range("H22").ClearFormats()
$ws.Range("H23").Value = "'1607409.56"
$ws.Range("H23").ClearFormats()
$ws.Range("H24").Value = "'1331915.13"
$ws.Range("H24").ClearFormats()
$ws.Range("H25").Value = "'57698.40"
$ws.Range("H25").ClearFormats()
$ws.Range("H26").Value = "'18770.50"
$ws.Range("H26").ClearFormats()
$ws.Range("H27").Value = "'47458.52"
$ws.Range("H27").ClearFormats()
$ws.Range("H28").Value = "'2730.00"
$ws.Range("H28").ClearFormats()
$ws.Range("H29").Value = "'1376334.73"
$ws.Range("H29").ClearFormats()
$ws.Range("H30").Value = "'299833.01"
$ws.Range("H30").ClearFormats()
$ws.Range("H31").Value = "'13500.00"
$ws.Range("H31").ClearFormats()
$ws.Range("H32").Value = "'119423.00"
$ws.Range("H32").ClearFormats()
$ws.Range("H33").Value = "'15840.00"
$ws.Range("H33").ClearFormats()
$ws.Range("H34").Value = "'3449.39"
$ws.Range("H34").ClearFormats()
$ws.Range("H35").Value = "'91043.74"
$ws.Range("H35").ClearFormats()
$ws.Range("H36").Value = "'11752.00"
$ws.Range("H36").ClearFormats()
$ws.Range("H37").Value = "'833.46"
$ws.Range("H37").ClearFormats()
$ws.Range("H38").Value = "'104414.23"
$ws.Range("H38").ClearFormats()
$ws.Range("H39").Value = "'32500.00"
$ws.Range("H39").ClearFormats()
$ws.Range("H40").Value = "'22290.00"
$ws.Range("H40").ClearFormats()
$ws.Range("H41").Value = "'10900.00"
$ws.Range("H41").ClearFormats()
$ws.Range("H42").Value = "'11700.00"
$ws.Range("H42").ClearFormats()
$ws.Range("H43").Value = "'160480.00"
$ws.Range("H43").ClearFormats()
$ws.Range("H44").Value = "'939.38"
$ws.Range("H44").ClearFormats()
$ws.Range("H45").Value = "'11540.00"
$ws.Range("H45").ClearFormats()
$ws.Range("H46").Value = "'443228.20"
$ws.Range("H46").ClearFormats()
$ws.Range("H47").Value = "'39348.81"
$ws.Range("H47").ClearFormats()
$ws.Range("H48").Value = "'466688.17"
$ws.Range("H48").ClearFormats()
$ws.Range("H49").Value = "'99645.00"
$ws.Range("H49").ClearFormats()
$ws.Range("H50").Value = "'82099.71"
$ws.Range("H50").ClearFormats()
$ws.Range("H51").Value = "'115000.00"
$ws.Range("H51").ClearFormats()
$ws.Range("H52").Value = "'21570.00"
$ws.Range("H52").ClearFormats()
$ws.Range("H53").Value = "'3000.00"
$ws.Range("H53").ClearFormats()
$ws.Range("H54").Value = "'2260.00"
$ws.Range("H54").ClearFormats()
$ws.Range("H56").Value = "'852322.00"
$ws.Range("H56").ClearFormats()
$ws.Range("H57").Value = "'314650.00"
$ws.Range("H57").ClearFormats()
$ws.Range("H58").Value = "'83945.00"
$ws.Range("H58").ClearFormats()
$ws.Range("H59").Value = "'40060.00"
$ws.Range("H59").ClearFormats()
$ws.Range("H60").Value = "'1576.60"
$ws.Range("H60").ClearFormats()
$ws.Range("H61").Value = "'67670.00"
$ws.Range("H61").ClearFormats()
$ws.Range("H62").Value = "'1899.99"
$ws.Range("H62").ClearFormats()
$ws.Range("H63").Value = "'4990.00"
$ws.Range("H63").ClearFormats()
$ws.Range("H64").Value = "'295279.89"
$ws.Range("H64").ClearFormats()
$ws.Range("H65").Value = "'192537.66"
$ws.Range("H65").ClearFormats()
$ws.Range("H66").Value = "'19803.70"
$ws.Range("H66").ClearFormats()
$ws.Range("H67").Value = "'8400.00"
$ws.Range("H67").ClearFormats()
$ws.Range("H68").Value = "'685.83"
$ws.Range("H68").ClearFormats()
$ws.Range("H69").Value = "'200.00"
$ws.Range("H69").ClearFormats()
$ws.Range("H70").Value = "'1150.00"
$ws.Range("H70").ClearFormats()
$ws.Range("H71").Value = "'10412.50"
$ws.Range("H71").ClearFormats()
$ws.Range("H72").Value = "'12166.88"
$ws.Range("H72").ClearFormats()
$ws.Range("H73").Value = "'6467.40"
$ws.Range("H73").ClearFormats()
$ws.Range("H74").Value = "'15941.84"
$ws.Range("H74").ClearFormats()
$ws.Range("H75").Value = "'57.38"
$ws.Range("H75").ClearFormats()
$ws.Range("H76").Value = "'7014.00"
$ws.Range("H76").ClearFormats()
$ws.Range("H77").Value = "'23260.00"
$ws.Range("H77").ClearFormats()
$ws.Range("H78").Value = "'48523.19"
$ws.Range("H78").ClearFormats()
$ws.Range("H79").Value = "'1000.00"
$ws.Range("H79").ClearFormats()
$ws.Range("H80").Value = "'2200.00"
$ws.Range("H80").ClearFormats()
$ws.Range("H81").Value = "'2250.00"
$ws.Range("H81").ClearFormats()
$ws.Range("H82").Value = "'1750.00"
$ws.Range("H82").ClearFormats()
$ws.Range("H83").Value = "'17181.96"
$ws.Range("H83").ClearFormats()
$ws.Range("H84").Value = "'1600.00"
$ws.Range("H84").ClearFormats()
$ws.Range("H85").Value = "'1317.99"
$ws.Range("H85").ClearFormats()
$ws.Range("H86").Value = "'232450.69"
$ws.Range("H86").ClearFormats()
$ws.Range("H87").Value = "'140.00"
$ws.Range("H87").ClearFormats()
$ws.Range("H88").Value = "'20448.63"
$ws.Range("H88").ClearFormats()
$ws.Range("H89").Value = "'18720.00"
$ws.Range("H89").ClearFormats()
$ws.Range("H90").Value = "'19200.00"
$ws.Range("H90").ClearFormats()
$ws.Range("H91").Value = "'1050.00"
$ws.Range("H91").ClearFormats()
$ws.Range("H92").Value = "'349650.00"
$ws.Range("H92").ClearFormats()
$ws.Range("H93").Value = "'18500.00"
$ws.Range("H93").ClearFormats()
$ws.Range("H94").Value = "'96000.00"
$ws.Range("H94").ClearFormats()
$ws.Range("H95").Value = "'26000.00"
$ws.Range("H95").ClearFormats()
$ws.Range("H96").Value = "'6825.00"
$ws.Range("H96").ClearFormats()
$ws.Range("H97").Value = "'20350.00"
$ws.Range("H97").ClearFormats()
$ws.Range("H98").Value = "'43800.00"
$ws.Range("H98").ClearFormats()
$ws.Range("H99").Value = "'14334.00"
$ws.Range("H99").ClearFormats()
$ws.Range("H100").Value = "'904.12"
$ws.Range("H100").ClearFormats()
$ws.Range("H101").Value = "'10890.07"
$ws.Range("H101").ClearFormats()
$ws.Range("H102").Value = "'2532.00"
$ws.Range("H102").ClearFormats()
$ws.Range("H103").Value = "'471.86"
$ws.Range("H103").ClearFormats()
$ws.Range("H104").Value = "'1070.00"
$ws.Range("H104").ClearFormats()
$ws.Range("H105").Value = "'21800.00"
$ws.Range("H105").ClearFormats()
$ws.Range("H106").Value = "'45417.40"
$ws.Range("H106").ClearFormats()
$ws.Range("H107").Value = "'3470.00"
$ws.Range("H107").ClearFormats()
$ws.Range("H108").Value = "'58798.88"
$ws.Range("H108").ClearFormats()
$ws.Range("H109").Value = "'406.90"
$ws.Range("H109").ClearFormats()
$ws.Range("H110").Value = "'4130.00"
$ws.Range("H110").ClearFormats()
$ws.Range("H111").Value = "'220255.00"
$ws.Range("H111").ClearFormats()
$ws.Range("H112").Value = "'10560.00"
$ws.Range("H112").ClearFormats()
$ws.Range("H113").Value = "'40700.00"
$ws.Range("H113").ClearFormats()
$ws.Range("H114").Value = "'129180.00"
$ws.Range("H114").ClearFormats()
$ws.Range("H115").Value = "'21026.04"
$ws.Range("H115").ClearFormats()
$ws.Range("H117").Value = "'6335.00"
$ws.Range("H117").ClearFormats()
$ws.Range("H118").Value = "'19650.34"
$ws.Range("H118").ClearFormats()
$ws.Range("H119").Value = "'500.78"
$ws.Range("H119").ClearFormats()
$ws.Range("H120").Value = "'40192.98"
$ws.Range("H120").ClearFormats()
$ws.Range("H121").Value = "'780.00"
$ws.Range("H121").ClearFormats()
$ws.Range("H122").Value = "'148304.00"
$ws.Range("H122").ClearFormats()
$ws.Range("H123").Value = "'5999.00"
$ws.Range("H123").ClearFormats()
$ws.Range("H124").Value = "'755.00"
$ws.Range("H124").ClearFormats()
$ws.Range("H125").Value = "'16100.00"
$ws.Range("H125").ClearFormats()
$ws.Range("H126").Value = "'499.00"
$ws.Range("H126").ClearFormats()
$ws.Range("H127").Value = "'5669.89"
$ws.Range("H127").ClearFormats()
$ws.Range("H128").Value = "'1020.00"
$ws.Range("H128").ClearFormats()
$ws.Range("H129").Value = "'5235.76"
$ws.Range("H129").ClearFormats()
$ws.Range("H130").Value = "'4972.00"
$ws.Range("H130").ClearFormats()
$ws.Range("H131").Value = "'13226.00"
$ws.Range("H131").ClearFormats()
$ws.Range("H132").Value = "'11478.19"
$ws.Range("H132").ClearFormats()
$ws.Range("H133").Value = "'32416.08"
$ws.Range("H133").ClearFormats()
$ws.Range("H134").Value = "'8407.20"
$ws.Range("H134").ClearFormats()
$ws.Range("H135").Value = "'2588.00"
$ws.Range("H135").ClearFormats()
$ws.Range("H136").Value = "'12930.00"
$ws.Range("H136").ClearFormats()
$ws.Range("H137").Value = "'638.80"
$ws.Range("H137").ClearFormats()
$ws.Range("H138").Value = "'2690.00"
$ws.Range("H138").ClearFormats()
$ws.Range("H139").Value = "'264.00"
$ws.Range("H139").ClearFormats()
$ws.Range("H140").Value = "'3027.00"
$ws.Range("H140").ClearFormats()
$ws.Range("H141").Value = "'13655.00"
$ws.Range("H141").ClearFormats()
$ws.Range("H142").Value = "'193.12"
$ws.Range("H142").ClearFormats()
$ws.Range("H144").Value = "'2993.40"
$ws.Range("H144").ClearFormats()
$ws.Range("H145").Value = "'5070.00"
$ws.Range("H145").ClearFormats()
$ws.Range("H146").Value = "'6750.00"
$ws.Range("H146").ClearFormats()
$ws.Range("H147").Value = "'10720.00"
$ws.Range("H147").ClearFormats()
$ws.Range("H148").Value = "'31000.00"
$ws.Range("H148").ClearFormats()
$ws.Range("H149").Value = "'11400.00"
$ws.Range("H149").ClearFormats()
$ws.Range("H150").Value = "'23238.00"
$ws.Range("H150").ClearFormats()
$ws.Range("H151").Value = "'32768.40"
$ws.Range("H151").ClearFormats()
$ws.Range("H152").Value = "'9030.00"
$ws.Range("H152").ClearFormats()
$ws.Range("H153").Value = "'13690.00"
$ws.Range("H153").ClearFormats()
$ws.Range("H154").Value = "'59146.52"
$ws.Range("H154").ClearFormats()
$ws.Range("H155").Value = "'3300.00"
$ws.Range("H155").ClearFormats()
$ws.Range("H159").Value = "'12000.00"
$ws.Range("H159").ClearFormats()
$ws.Range("H160").Value = "'25000.00"
$ws.Range("H160").ClearFormats()
$ws.Range("H161").Value = "'58000.00"
$ws.Range("H161").ClearFormats()
$ws.Range("H164").Value = "'138000.00"
$ws.Range("H164").ClearFormats()
$ws.Range("H165").Value = "'50000.00"
$ws.Range("H165").ClearFormats()
$ws.Range("H167").Value = "'21009.57"
$ws.Range("H167").ClearFormats()
$ws.Range("H168").Value = "'7780.08"
$ws.Range("H168").ClearFormats()
$ws.Range("H169").Value = "'47544.99"
$ws.Range("H169").ClearFormats()
$ws.Range("H170").Value = "'8610.00"
$ws.Range("H170").ClearFormats()
$ws.Range("H171").Value = "'10390.00"
$ws.Range("H171").ClearFormats()
$ws.Range("H172").Value = "'20726.85"
$ws.Range("H172").ClearFormats()
$ws.Range("H173").Value = "'2744400.00"
$ws.Range("H173").ClearFormats()
$ws.Range("H174").Value = "'70925.00"
$ws.Range("H174").ClearFormats()
$ws.Range("H176").Value = "'33000.00"
$ws.Range("H176").ClearFormats()
$ws.Range("H178").Value = "'66000.00"
$ws.Range("H178").ClearFormats()
$ws.Range("H182").Value = "'30000.00"
$ws.Range("H182").ClearFormats()
$ws.Range("H183").Value = "'17000.00"
$ws.Range("H183").ClearFormats()
$ws.Range("H188").Value = "'18000.00"
$ws.Range("H188").ClearFormats()
$ws.Range("H189").Value = "'11500.00"
$ws.Range("H189").ClearFormats()
$ws.Range("H193").Value = "'42000.00"
$ws.Range("H193").ClearFormats()
$ws.Range("H197").Value = "'10000.00"
$ws.Range("H197").ClearFormats()
$ws.Range("H198").Value = "'188166.11"
$ws.Range("H198").ClearFormats()
$ws.Range("H199").Value = "'32000.00"
$ws.Range("H199").ClearFormats()
$ws.Range("H201").Value = "'24000.00"
$ws.Range("H201").ClearFormats()
$ws.Range("H203").Value = "'25300.00"
$ws.Range("H203").ClearFormats()
$ws.Range("H204").Value = "'14000.00"
$ws.Range("H204").ClearFormats()
$ws.Range("H207").Value = "'82600.00"
$ws.Range("H207").ClearFormats()
$ws.Range("H208").Value = "'7000.00"
$ws.Range("H208").ClearFormats()
$ws.Range("H209").Value = "'82960.00"
$ws.Range("H209").ClearFormats()
$ws.Range("H210").Value = "'20000.00"
$ws.Range("H210").ClearFormats()
$ws.Range("H211").Value = "'9500.00"
$ws.Range("H211").ClearFormats()
$ws.Range("H212").Value = "'9800.00"
$ws.Range("H212").ClearFormats()
$ws.Range("H213").Value = "'27000.00"
$ws.Range("H213").ClearFormats()
$ws.Range("H214").Value = "'2700.00"
$ws.Range("H214").ClearFormats()
$ws.Range("H216").Value = "'46730.00"
$ws.Range("H216").ClearFormats()
$ws.Range("H217").Value = "'9000.00"
$ws.Range("H217").ClearFormats()
$ws.Range("H218").Value = "'290.84"
$ws.Range("H218").ClearFormats()
$ws.Range("H219").Value = "'10610.00"
$ws.Range("H219").ClearFormats()
$ws.Range("H220").Value = "'38000.00"
$ws.Range("H220").ClearFormats()
$ws.Range("H221").Value = "'362.00"
$ws.Range("H221").ClearFormats()
$ws.Range("H222").Value = "'32600.00"
$ws.Range("H222").ClearFormats()
$ws.Range("H223").Value = "'43765.00"
$ws.Range("H223").ClearFormats()
$ws.Range("H224").Value = "'9546.62"
$ws.Range("H224").ClearFormats()
$ws.Range("H225").Value = "'29.30"
$ws.Range("H225").ClearFormats()
$ws.Range("H226").Value = "'106800.00"
$ws.Range("H226").ClearFormats()
$ws.Range("H227").Value = "'2380.00"
$ws.Range("H227").ClearFormats()
$ws.Range("H228").Value = "'4465.00"
$ws.Range("H228").ClearFormats()
$ws.Range("H229").Value = "'2400.00"
$ws.Range("H229").ClearFormats()
$ws.Range("H230").Value = "'600.00"
$ws.Range("H230").ClearFormats()
$ws.Range("H231").Value = "'16800.44"
$ws.Range("H231").ClearFormats()
$ws.Range("H232").Value = "'6000.00"
$ws.Range("H232").ClearFormats()
$ws.Range("H233").Value = "'16905.45"
$ws.Range("H233").ClearFormats()
$ws.Range("H234").Value = "'16404.00"
$ws.Range("H234").ClearFormats()
$ws.Range("H235").Value = "'4382.82"
$ws.Range("H235").ClearFormats()
$ws.Range("H236").Value = "'44450.00"
$ws.Range("H236").ClearFormats()
$ws.Range("H237").Value = "'3100.00"
$ws.Range("H237").ClearFormats()
$ws.Range("H238").Value = "'3040.00"
$ws.Range("H238").ClearFormats()
$ws.Range("H239").Value = "'3660.19"
$ws.Range("H239").ClearFormats()
$ws.Range("H240").Value = "'7071.00"
$ws.Range("H240").ClearFormats()
$ws.Range("H241").Value = "'1214.52"
$ws.Range("H241").ClearFormats()
$ws.Range("H242").Value = "'1054.00"
$ws.Range("H242").ClearFormats()
$ws.Range("H243").Value = "'76160.00"
$ws.Range("H243").ClearFormats()
$ws.Range("H244").Value = "'16058.00"
$ws.Range("H244").ClearFormats()
$ws.Range("H245").Value = "'927.30"
$ws.Range("H245").ClearFormats()
$ws.Range("H246").Value = "'15940.00"
$ws.Range("H246").ClearFormats()
$ws.Range("H247").Value = "'12753.75"
$ws.Range("H247").ClearFormats()
$ws.Range("H253").Value = "'55000.00"
$ws.Range("H253").ClearFormats()
$ws.Range("H255").Value = "'40000.00"
$ws.Range("H255").ClearFormats()
$ws.Range("H257").Value = "'80000.00"
$ws.Range("H257").ClearFormats()
$ws.Range("H258").Value = "'8600.00"
$ws.Range("H258").ClearFormats()
$ws.Range("H259").Value = "'7500.00"
$ws.Range("H259").ClearFormats()
$ws.Range("H260").Value = "'4800.00"
$ws.Range("H260").ClearFormats()
$ws.Range("H261").Value = "'7281492.38"
$ws.Range("H261").ClearFormats()
$ws.Range("H262").Value = "'1838080.00"
$ws.Range("H262").ClearFormats()
$ws.Range("H263").Value = "'48120.00"
$ws.Range("H263").ClearFormats()
$ws.Range("H264").Value = "'8000.00"
$ws.Range("H264").ClearFormats()
$ws.Range("H265").Value = "'9693700.00"
$ws.Range("H265").ClearFormats()
$ws.Range("H266").Value = "'92389.78"
$ws.Range("H266").ClearFormats()
$ws.Range("H268").Value = "'282300.00"
$ws.Range("H268").ClearFormats()
$ws.Range("H269").Value = "'253000.00"
$ws.Range("H269").ClearFormats()
$ws.Range("H270").Value = "'258800.00"
$ws.Range("H270").ClearFormats()
$ws.Range("H271").Value = "'251600.00"
$ws.Range("H271").ClearFormats()
$ws.Range("H273").Value = "'480350.00"
$ws.Range("H273").ClearFormats()
$ws.Range("H275").Value = "'603450.00"
$ws.Range("H275").ClearFormats()
$ws.Range("H276").Value = "'561000.00"
$ws.Range("H276").ClearFormats()
$ws.Range("H277").Value = "'331900.00"
$ws.Range("H277").ClearFormats()
$ws.Range("H278").Value = "'250000.00"
$ws.Range("H278").ClearFormats()
$ws.Range("H280").Value = "'398000.00"
$ws.Range("H280").ClearFormats()
$ws.Range("H281").Value = "'493200.00"
$ws.Range("H281").ClearFormats()
$ws.Range("H282").Value = "'714300.00"
$ws.Range("H282").ClearFormats()
$ws.Range("H283").Value = "'470000.00"
$ws.Range("H283").ClearFormats()
$ws.Range("H284").Value = "'714900.00"
$ws.Range("H284").ClearFormats()
$ws.Range("H285").Value = "'500000.00"
$ws.Range("H285").ClearFormats()
$ws.Range("H286").Value = "'257750.00"
$ws.Range("H286").ClearFormats()
$ws.Range("H287").Value = "'14500.00"
$ws.Range("H287").ClearFormats()
$ws.Range("H288").Value = "'28000.00"
$ws.Range("H288").ClearFormats()
$ws.Range("H289").Value = "'5203738.70"
$ws.Range("H289").ClearFormats()
$ws.Range("H290").Value = "'2843700.56"
$ws.Range("H290").ClearFormats()
$ws.Range("H291").Value = "'1156.17"
$ws.Range("H291").ClearFormats()
$ws.Range("H292").Value = "'56200.00"
$ws.Range("H292").ClearFormats()
$ws.Range("H293").Value = "'112200.00"
$ws.Range("H293").ClearFormats()
$ws.Range("H294").Value = "'186000.00"
$ws.Range("H294").ClearFormats()
$ws.Range("H295").Value = "'15000.00"
$ws.Range("H295").ClearFormats()
$ws.Range("H296").Value = "'10500.00"
$ws.Range("H296").ClearFormats()
$ws.Range("H297").Value = "'4200.00"
$ws.Range("H297").ClearFormats()
$ws.Range("H298").Value = "'990000.00"
$ws.Range("H298").ClearFormats()
$ws.Range("H299").Value = "'24880.00"
$ws.Range("H299").ClearFormats()
$ws.Range("H300").Value = "'58400.00"
$ws.Range("H300").ClearFormats()
$ws.Range("H301").Value = "'98000.00"
$ws.Range("H301").ClearFormats()
$ws.Range("H302").Value = "'133450.00"
$ws.Range("H302").ClearFormats()
$ws.Range("H303").Value = "'343500.00"
$ws.Range("H303").ClearFormats()
$ws.Range("H304").Value = "'6676.00"
$ws.Range("H304").ClearFormats()
$ws.Range("H305").Value = "'4000.00"
$ws.Range("H305").ClearFormats()
$ws.Range("H306").Value = "'14700.00"
$ws.Range("H306").ClearFormats()
$ws.Range("H307").Value = "'28696.00"
$ws.Range("H307").ClearFormats()
